# NewNutrientData.xlsx revision:
#  - Drop the "Distance" (G) and "Notes" (H) columns entirely — the field
#    notes captured during sampling are no longer part of the published
#    data table.
#  - Consolidate the "Pre-Amendment" / "Pre-Treatment" timepoint labels
#    into the existing "1hr" timepoint (they described the same baseline
#    draw, just inconsistently labelled across sampling days).
#  - Give the new first/last column ("SampleID" header) the bold header
#    styling that used to live on the removed "Notes" header.
#  - Set an explicit width on the now-last column (Measurement).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove Notes (H) then Distance (G) — delete from the right so column
# indices for the earlier deletion stay valid.
$ws.Columns.Item(8).Delete()
$ws.Columns.Item(7).Delete()

# Normalize the Timepoint column: both stand-in baseline labels become "1hr".
$timepointRange = $ws.Range("D1:D121")
$timepointRange.Replace("Pre-Amendment", "1hr")
$timepointRange.Replace("Pre-Treatment", "1hr")

# Explicit width for column F (Measurement), now the last column.
$ws.Columns.Item(6).ColumnWidth = 11.14

# Bold header styling on A1 (SampleID), matching the old Notes header look.
$ws.Range("A1").Font.Bold = $true
